# Apply cryptocurrency price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.925.13'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '3.638.22'
$ws.Range('E3').Value = '  +3.96%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''605.20'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '''200.41'
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '''0.219'
$ws.Range('E9').Value = '  +9.65%  '
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').Value = '''53.93'
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').Value = '''0.0000307'
$ws.Range('E12').Value = '  +2.38%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '4.208.82'
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('D15').Value = '''636.04'
$ws.Range('E15').Value = '  +7.09%  '
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').Value = '70.934.30'
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').Value = '3.623.75'
$ws.Range('E18').Value = '  +3.76%  '
$ws.Range('D19').Value = '''19.08'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').Value = '''18.28'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').Value = '''5.39'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').Value = '''104.20'
$ws.Range('E24').Value = '  +2.10%  '
$ws.Range('D25').Value = '''4.64'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  -4.48%  '
$ws.Range('D27').Value = '''10.56'
$ws.Range('E27').Value = '  -2.69%  '
$ws.Range('D28').Value = '''9.76'
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('D29').Value = '''33.80'
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('D30').Value = '''4.80'
$ws.Range('E30').Value = '  +13.34%  '
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').Value = '''12.28'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('D34').Value = '''63.41'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = '0.0₃0882'
$ws.Range('E35').Value = '  +6.28%  '
$ws.Range('D36').Value = '3.995.15'
$ws.Range('E36').Value = '  +7.48%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = '''517.02'
$ws.Range('E38').Value = '  +8.56%  '
$ws.Range('D39').Value = '''3.04'
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '''0.390'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '''36.77'
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('D44').Value = '''0.0463'
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('D45').Value = '''3.51'
$ws.Range('E45').Value = '  +7.25%  '
$ws.Range('D46').Value = '''2.92'
$ws.Range('E46').Value = '  +4.09%  '
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').Value = '''8.64'
$ws.Range('E48').Value = '  +2.67%  '
$ws.Range('E49').Value = '  -0.36%  '
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('E51').Value = '  +1.28%  '
